$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws3 = $wb.Worksheets.Item("model")
$ws5 = $wb.Worksheets.Item("choices")

# ---------- survey sheet: insert a new row 22 (assign / is_active / TRUE) ----------
# shifting the old row 22 (note / "Press "Next" to finish.") down to row 23.
# Row 21 is cleared first so the row-insert doesn't inherit/duplicate its
# formatting into stray cells, then row 21 is restored from the saved values.
$a21 = $ws1.Range("A21").Text
$b21 = $ws1.Range("B21").Text
$c21 = $ws1.Range("C21").Text
$d21 = $ws1.Range("D21").Text
$ws1.Rows.Item(21).Clear()
$ws1.Rows.Item(22).Insert()
$ws1.Range("A21").Value = $a21
$ws1.Range("B21").Value = $b21
$ws1.Range("C21").Value = $c21
$ws1.Range("D21").Value = $d21

$ws1.Range("A22").Value = "assign"
$ws1.Range("B22").Value = "is_active"
$ws1.Range("E22").Value = $true

# ---------- model sheet: add the matching field-type rows ----------
$ws3.Range("B21").Value = "is_active"
$ws3.Range("B22").Value = "disabled_reason"
$ws3.Range("A21").Value = "boolean"
$ws3.Range("A23").Value = "boolean"
$ws3.Range("B23").Value = "is_override"
$ws3.Range("A22").Value = "string"
$ws3.Range("B25").Value = " "

# ---------- selections on each sheet ----------
$ws3.Range("B25").Select() | Out-Null
$ws5.Range("E11").Select() | Out-Null
$ws1.Range("F22").Select() | Out-Null

# ---------- make "survey" the active tab ----------
$ws1.Activate() | Out-Null
